$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NATMI TPM recalculation: drop the three "Target cluster = ECs" rows
# (old rows 2, 5, 8), shifting the remaining rows up, then refresh the
# NATMI-derived numeric columns for the surviving Ngf-Ntrk1 pairs.
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(2).Delete()

# Row 2
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.322531
$ws.Cells.Item(2, 8).Value = 0.967593
$ws.Cells.Item(2, 9).Value = 0.01892149513432853
$ws.Cells.Item(2, 10).Value = 0.01892149513432853
$ws.Cells.Item(2, 15).Value = 0.8034178059852001
$ws.Cells.Item(2, 16).Value = 0.8034178059852
$ws.Cells.Item(2, 17).Value = 0.02886598694833334
$ws.Cells.Item(2, 18).Value = 0.259793882535
$ws.Cells.Item(2, 19).Value = 0.01520186610678187
$ws.Cells.Item(2, 20).Value = 0.01520186610678186

# Row 3
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.322531
$ws.Cells.Item(3, 8).Value = 0.967593
$ws.Cells.Item(3, 9).Value = 0.01892149513432853
$ws.Cells.Item(3, 10).Value = 0.01892149513432853
$ws.Cells.Item(3, 15).Value = 0.1965821940147999
$ws.Cells.Item(3, 16).Value = 0.1965821940147999
$ws.Cells.Item(3, 17).Value = 0.007062998858666667
$ws.Cells.Item(3, 18).Value = 0.06356698972800001
$ws.Cells.Item(3, 19).Value = 0.003719629027546665
$ws.Cells.Item(3, 20).Value = 0.003719629027546664

# Row 4
$ws.Cells.Item(4, 9).Value = 0.0261208867009986
$ws.Cells.Item(4, 10).Value = 0.0261208867009986
$ws.Cells.Item(4, 15).Value = 0.8034178059852001
$ws.Cells.Item(4, 16).Value = 0.8034178059852
$ws.Cells.Item(4, 19).Value = 0.02098598548370428
$ws.Cells.Item(4, 20).Value = 0.02098598548370428

# Row 5
$ws.Cells.Item(5, 9).Value = 0.0261208867009986
$ws.Cells.Item(5, 10).Value = 0.0261208867009986
$ws.Cells.Item(5, 15).Value = 0.1965821940147999
$ws.Cells.Item(5, 16).Value = 0.1965821940147999
$ws.Cells.Item(5, 19).Value = 0.005134901217294313
$ws.Cells.Item(5, 20).Value = 0.005134901217294313

# Row 6
$ws.Cells.Item(6, 9).Value = 0.954957618164673
$ws.Cells.Item(6, 10).Value = 0.954957618164673
$ws.Cells.Item(6, 15).Value = 0.8034178059852001
$ws.Cells.Item(6, 16).Value = 0.8034178059852
$ws.Cells.Item(6, 19).Value = 0.767229954394714
$ws.Cells.Item(6, 20).Value = 0.7672299543947139

# Row 7
$ws.Cells.Item(7, 9).Value = 0.954957618164673
$ws.Cells.Item(7, 10).Value = 0.954957618164673
$ws.Cells.Item(7, 15).Value = 0.1965821940147999
$ws.Cells.Item(7, 16).Value = 0.1965821940147999
$ws.Cells.Item(7, 19).Value = 0.187727663769959
$ws.Cells.Item(7, 20).Value = 0.187727663769959
